$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 16 data, following the same pattern as row 15 (HexGrid-60degTilt5degRes)

# Copy formatting from A15 (bold/centered/bordered style) to A16
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.96339791054673
$ws.Range("D16").Value = 1.224158669892441
$ws.Range("E16").Value = 0.936824953861738
$ws.Range("F16").Value = 0.96339791054673
$ws.Range("G16").Value = 1.11310358006911
$ws.Range("H16").Value = 0.8511507835846838
$ws.Range("I16").Value = 0.9392363912951788
$ws.Range("J16").Value = 1.224158669892441
$ws.Range("K16").Value = 1.080491811877089
$ws.Range("L16").Value = 1.02194486121191
$ws.Range("M16").Value = 1.004645381541647
